# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit:
#   1) Slide 6's table switches from the deck's custom "Table_0" style to the
#      built-in table style {49C8AEAC-D410-4C97-B592-2A2E9B0C92CC}.
#   2) The presentation's applied colour theme changes from the "Integral"
#      palette back to the stock "Office" palette (dk1/lt1 are identical in
#      both palettes, so only the other 10 theme colours actually move).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table on slide 6 (title + table are the two shapes; table is #2) ---
$slide = $p.Slides.Item(6)
$table = $slide.Shapes.Item(2).Table
$table.ApplyStyle("{49C8AEAC-D410-4C97-B592-2A2E9B0C92CC}")

# --- 2) Swap the theme's colour scheme from Integral back to Office -------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
